# Apply the update described by the commit: refresh odds data causing
# several existing rows to swap their home/away contents (columns F:V),
# and append three brand-new match rows (107-109) at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param([int]$Row1, [int]$Row2)

    $range1 = $ws.Range("F" + $Row1 + ":V" + $Row1)
    $range2 = $ws.Range("F" + $Row2 + ":V" + $Row2)

    $v1 = $range1.Value2
    $v2 = $range2.Value2

    $range1.Value = $v2
    $range2.Value = $v1
}

# Simple pairwise swaps of match content (F:V) between rows.
Swap-Rows 2 3
Swap-Rows 4 5
Swap-Rows 26 27
Swap-Rows 53 55
Swap-Rows 69 70
Swap-Rows 86 87

# Three-way rotation among rows 56, 57, 58:
#   new 56 <- old 58, new 57 <- old 56, new 58 <- old 57
$v56 = $ws.Range("F56:V56").Value2
$v57 = $ws.Range("F57:V57").Value2
$v58 = $ws.Range("F58:V58").Value2

$ws.Range("F56:V56").Value = $v58
$ws.Range("F57:V57").Value = $v56
$ws.Range("F58:V58").Value = $v57

# Append three new rows (107-109) at the bottom, copying formatting from
# the existing last row (106) and then filling in the new values.
$ws.Range("A106:V106").Copy()
$ws.Range("A107:V109").PasteSpecial(-4122)

# Row 107
$ws.Cells.Item(107, 1).Value = 106
$ws.Cells.Item(107, 2).Value = "italy"
$ws.Cells.Item(107, 3).Value = "serie-a"
$ws.Cells.Item(107, 4).Value = "2023-2024"
$ws.Cells.Item(107, 5).Value = 45235.625
$ws.Cells.Item(107, 6).Value = "Cagliari"
$ws.Cells.Item(107, 7).Value = 2
$ws.Cells.Item(107, 8).Value = "Genoa"
$ws.Cells.Item(107, 9).Value = 1
$ws.Cells.Item(107, 10).Value = 2.6
$ws.Cells.Item(107, 11).Value = "23/10/2023 15:49"
$ws.Cells.Item(107, 12).Value = 2.76
$ws.Cells.Item(107, 13).Value = "05/11/2023 14:58"
$ws.Cells.Item(107, 14).Value = 3.21
$ws.Cells.Item(107, 15).Value = "23/10/2023 15:49"
$ws.Cells.Item(107, 16).Value = 2.98
$ws.Cells.Item(107, 17).Value = "05/11/2023 14:59"
$ws.Cells.Item(107, 18).Value = 2.79
$ws.Cells.Item(107, 19).Value = "23/10/2023 15:49"
$ws.Cells.Item(107, 20).Value = 3.05
$ws.Cells.Item(107, 21).Value = "05/11/2023 14:59"
$ws.Cells.Item(107, 22).Value = "https://www.betexplorer.com/football/italy/serie-a/cagliari-genoa/MN8LgIJo/"

# Row 108
$ws.Cells.Item(108, 1).Value = 107
$ws.Cells.Item(108, 2).Value = "italy"
$ws.Cells.Item(108, 3).Value = "serie-a"
$ws.Cells.Item(108, 4).Value = "2023-2024"
$ws.Cells.Item(108, 5).Value = 45235.75
$ws.Cells.Item(108, 6).Value = "AS Roma"
$ws.Cells.Item(108, 7).Value = 2
$ws.Cells.Item(108, 8).Value = "Lecce"
$ws.Cells.Item(108, 9).Value = 1
$ws.Cells.Item(108, 10).Value = 1.49
$ws.Cells.Item(108, 11).Value = "22/10/2023 12:02"
$ws.Cells.Item(108, 12).Value = 1.56
$ws.Cells.Item(108, 13).Value = "05/11/2023 17:58"
$ws.Cells.Item(108, 14).Value = 4.12
$ws.Cells.Item(108, 15).Value = "22/10/2023 12:02"
$ws.Cells.Item(108, 16).Value = 4.01
$ws.Cells.Item(108, 17).Value = "05/11/2023 17:59"
$ws.Cells.Item(108, 18).Value = 6.82
$ws.Cells.Item(108, 19).Value = "22/10/2023 12:02"
$ws.Cells.Item(108, 20).Value = 7.2
$ws.Cells.Item(108, 21).Value = "05/11/2023 17:59"
$ws.Cells.Item(108, 22).Value = "https://www.betexplorer.com/football/italy/serie-a/as-roma-lecce/G4AyjzJA/"

# Row 109
$ws.Cells.Item(109, 1).Value = 108
$ws.Cells.Item(109, 2).Value = "italy"
$ws.Cells.Item(109, 3).Value = "serie-a"
$ws.Cells.Item(109, 4).Value = "2023-2024"
$ws.Cells.Item(109, 5).Value = 45235.86458333334
$ws.Cells.Item(109, 6).Value = "Fiorentina"
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = "Juventus"
$ws.Cells.Item(109, 9).Value = 1
$ws.Cells.Item(109, 10).Value = 2.77
$ws.Cells.Item(109, 11).Value = "22/10/2023 12:02"
$ws.Cells.Item(109, 12).Value = 2.9
$ws.Cells.Item(109, 13).Value = "05/11/2023 20:44"
$ws.Cells.Item(109, 14).Value = 3.21
$ws.Cells.Item(109, 15).Value = "22/10/2023 12:02"
$ws.Cells.Item(109, 16).Value = 3.15
$ws.Cells.Item(109, 17).Value = "05/11/2023 20:43"
$ws.Cells.Item(109, 18).Value = 2.77
$ws.Cells.Item(109, 19).Value = "22/10/2023 12:02"
$ws.Cells.Item(109, 20).Value = 2.75
$ws.Cells.Item(109, 21).Value = "05/11/2023 20:44"
$ws.Cells.Item(109, 22).Value = "https://www.betexplorer.com/football/italy/serie-a/fiorentina-juventus/0E7PhxZi/"
